$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (personas-mano-obra-familiar-con-remuneracion) was re-curated
# from a dimension to a measure, and its mapping file is no longer used.
$ws.Range("G2").Value = "iaest-measure:personas-mano-obra-familiar-con-remuneracion"
$ws.Range("G3").Value = "medida"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("G5").Clear()

# Column J (provincia) is re-curated from sdmx-dimension:refArea /
# URI-Provincia into a plain iaest measure.
$ws.Range("J2").Value = "iaest-measure:provincia"
$ws.Range("J3").Value = "medida"
$ws.Range("J4").Value = "xsd:int"
